# Insert a new data row at row 623 (pushing the existing rows 623:750 down to
# 624:751) and populate it with the new record. This mirrors the source
# workbook's "Fruta / hortaliza, semanal" commit, which adds one new weekly
# price observation near the top of the historical block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 623:750 down one row, duplicating row 622's formatting
# (same as Excel does when you right-click a row header -> Insert).
$ws.Rows("623:623").Insert()

# Populate the newly-inserted row with the new record's values.
$ws.Cells.Item(623, 1).Value = 3
$ws.Cells.Item(623, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(623, 3).Value = "Coquimbo"
$ws.Cells.Item(623, 4).Value = 45258
$ws.Cells.Item(623, 5).Value = 5
$ws.Cells.Item(623, 6).Value = "Fruta"
$ws.Cells.Item(623, 7).Value = 100108
$ws.Cells.Item(623, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(623, 9).Value = 100108002
$ws.Cells.Item(623, 10).Value = "Mango"
$ws.Cells.Item(623, 11).Value = "Sin especificar"
$ws.Cells.Item(623, 12).Value = "Primera"
$ws.Cells.Item(623, 13).Value = 228
$ws.Cells.Item(623, 14).Value = 11000
$ws.Cells.Item(623, 15).Value = 11000
$ws.Cells.Item(623, 16).Value = 11000
$ws.Cells.Item(623, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(623, 18).Value = "Perú"
$ws.Cells.Item(623, 19).Value = 2750
$ws.Cells.Item(623, 20).Value = 4
